$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('AK_1.png', 0.364, 0.001, 0.095, 0.367, 'Tidak Diketahui', 'Salah')
    ,@('AK_2.png', 0.681, 0.002, 0.337, 0.7, 'Akhlak Kamiswara', 'Benar')
    ,@('AK_3.png', 0.676, 0.002, 0.427, 0.733, 'Akhlak Kamiswara', 'Benar')
    ,@('AK_4.png', 0.339, 0.001, 0.43, 0.267, 'Tidak Diketahui', 'Salah')
    ,@('AK_5.png', 0.323, 0.001, 0.439, 0.333, 'Tidak Diketahui', 'Salah')
    ,@('MIB_1.png', 1.463, 0.005, 0.08599999999999999, 0.5, 'Muhammad Iqbal Baqi', 'Benar')
    ,@('MIB_2.png', 0.986, 0.003, 0.275, 0.633, 'Muhammad Iqbal Baqi', 'Benar')
    ,@('MIB_3.png', 1.621, 0.005, 0.54, 0.9, 'Muhammad Iqbal Baqi', 'Benar')
    ,@('MIB_4.png', 0.984, 0.003, 0.199, 0.7, 'Muhammad Iqbal Baqi', 'Benar')
    ,@('MIB_5.png', 1.395, 0.005, 0.216, 0.667, 'Muhammad Iqbal Baqi', 'Benar')
    ,@('AAH_1.png', 0.67, 0.002, 0.24, 0.9, 'Andrea Ayunove Hutami', 'Benar')
    ,@('AAH_2.png', 0.9399999999999999, 0.003, 0.6, 1, 'Andrea Ayunove Hutami', 'Benar')
    ,@('AAH_3.png', 0.733, 0.002, 0.244, 0.967, 'Andrea Ayunove Hutami', 'Benar')
    ,@('TI_1.png', 0.765, 0.003, 0.197, 0.533, 'Muhammad Iqbal Baqi', 'Salah')
    ,@('TI_2.png', 0.798, 0.003, 0.327, 0.533, 'Toni Ismail', 'Benar')
    ,@('TI_3.png', 0.553, 0.002, 0.554, 0.433, 'Tidak Diketahui', 'Salah')
    ,@('TI_4.png', 0.527, 0.002, 0.281, 0.6, 'Toni Ismail', 'Benar')
    ,@('TI_5.png', 0.927, 0.003, 0.356, 0.533, 'Toni Ismail', 'Benar')
    ,@('RAS_1.png', 0.478, 0.002, 0.371, 0.333, 'Tidak Diketahui', 'Salah')
    ,@('RAS_2.png', 0.862, 0.003, 0.269, 0.5, 'Ridha Ayu Salsabila', 'Benar')
    ,@('RAS_3.png', 0.356, 0.001, 0.343, 0.233, 'Tidak Diketahui', 'Salah')
    ,@('RAS_4.png', 1.062, 0.004, 0.097, 0.333, 'Tidak Diketahui', 'Salah')
    ,@('RAS_5.png', 1.003, 0.003, 0.345, 0.533, 'Ridha Ayu Salsabila', 'Benar')
    ,@('RR_1.png', 1.342, 0.004, 0.41, 0.667, 'Rafiqo Rapitasari', 'Benar')
    ,@('RR_2.png', 1.221, 0.004, 0.468, 0.5669999999999999, 'Rafiqo Rapitasari', 'Benar')
    ,@('RR_3.png', 0.924, 0.003, 0.096, 0.5669999999999999, 'Rafiqo Rapitasari', 'Benar')
    ,@('RR_4.png', 1.218, 0.004, 0.425, 0.633, 'Rafiqo Rapitasari', 'Benar')
    ,@('RR_5.png', 1.334, 0.004, 0.445, 0.733, 'Rafiqo Rapitasari', 'Benar')
    ,@('AR_1.png', 0.6879999999999999, 0.002, 0.317, 0.667, 'Arizli Romadhon', 'Benar')
    ,@('GA_1.png', 1.569, 0.005, 0.452, 0.9330000000000001, 'Gege Ardiyansyah', 'Benar')
    ,@('GA_2.png', 0.488, 0.002, 0.25, 0.7, 'Gege Ardiyansyah', 'Benar')
    ,@('GA_3.png', 0.576, 0.002, 0.169, 0.6, 'Gege Ardiyansyah', 'Benar')
    ,@('FY_1.png', 0.948, 0.003, 0.232, 0.6, 'Fanny Yusuf', 'Benar')
    ,@('FY_2.png', 1.602, 0.005, 0.47, 0.3, 'Tidak Diketahui', 'Salah')
    ,@('FY_3.png', 1.323, 0.004, 0.297, 0.467, 'Tidak Diketahui', 'Salah')
    ,@('FY_4.png', 1.128, 0.004, 0.227, 0.367, 'Tidak Diketahui', 'Salah')
    ,@('TO_1.png', 0.63, 0.002, 0.357, 0.533, 'Tiara Oktavian', 'Benar')
    ,@('TO_2.png', 0.899, 0.003, 0.368, 0.533, 'Tiara Oktavian', 'Benar')
    ,@('TO_3.png', 0.6870000000000001, 0.002, 0.391, 0.5, 'Tiara Oktavian', 'Benar')
    ,@('TO_4.png', 6.848, 0.023, 0.111, 0.533, 'Tiara Oktavian', 'Benar')
    ,@('TO_5.png', 6.746, 0.022, 0.08799999999999999, 0.667, 'Tiara Oktavian', 'Benar')
    ,@('TD_1.png', 3.557, 0.012, 0.038, 0.3, 'Tidak Diketahui', 'Benar')
    ,@('TD_2.png', 4.801, 0.016, 0.053, 0.367, 'Tidak Diketahui', 'Benar')
    ,@('TD_3.png', 0.931, 0.003, 0.172, 0.367, 'Tidak Diketahui', 'Benar')
    ,@('TD_4.png', 0.79, 0.003, 0.097, 0.2, 'Tidak Diketahui', 'Benar')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = [double]$row[1]
    $ws.Cells.Item($r, 3).Value = [double]$row[2]
    $ws.Cells.Item($r, 4).Value = [double]$row[3]
    $ws.Cells.Item($r, 5).Value = [double]$row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
